$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $v = $cell.Value2
        if ($v -ne $null -and $v -is [string] -and $v.Contains($oldVersion)) {
            $cell.Value2 = $v.Replace($oldVersion, $newVersion)
        }
    }
}
